$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 odds update ---
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67

# --- Row 3 odds update ---
$ws.Range("I3").Value = 3.7
$ws.Range("K3").Value = 2.05
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("AC3").Value = 8
$ws.Range("AP3").Value = 23
$ws.Range("AZ3").Value = 67
$ws.Range("BB3").Value = 251

# --- Insert a new row at position 6 (shifts old rows 6,7 down to 7,8) ---
$ws.Range("A6:BD6").EntireRow.Insert()

# --- Populate the newly inserted row 6 (Defensor Sp. vs Progreso match) ---
$ws.Cells.Item(6, 1).Value = "8hY9xlYF"
$ws.Cells.Item(6, 2).Value = "13/11/2024"
$ws.Cells.Item(6, 3).Value = "21:30"
$ws.Cells.Item(6, 4).Value = "URUGUAY - PRIMERA DIVISION"
$ws.Cells.Item(6, 5).Value = "Defensor Sp."
$ws.Cells.Item(6, 6).Value = "Progreso"
$ws.Cells.Item(6, 7).Value = 1.73
$ws.Cells.Item(6, 8).Value = 3.4
$ws.Cells.Item(6, 9).Value = 5.5
$ws.Cells.Item(6, 10).Value = 2.4
$ws.Cells.Item(6, 11).Value = 2.05
$ws.Cells.Item(6, 12).Value = 5.5
$ws.Cells.Item(6, 13).Value = 1.08
$ws.Cells.Item(6, 14).Value = 8
$ws.Cells.Item(6, 15).Value = 1.4
$ws.Cells.Item(6, 16).Value = 2.75
$ws.Cells.Item(6, 17).Value = 2.25
$ws.Cells.Item(6, 18).Value = 1.62
$ws.Cells.Item(6, 19).Value = 1.5
$ws.Cells.Item(6, 20).Value = 2.5
$ws.Cells.Item(6, 21).Value = 2.1
$ws.Cells.Item(6, 22).Value = 1.67
$ws.Cells.Item(6, 23).Value = 5.5
$ws.Cells.Item(6, 24).Value = 7
$ws.Cells.Item(6, 25).Value = 9
$ws.Cells.Item(6, 26).Value = 13
$ws.Cells.Item(6, 27).Value = 17
$ws.Cells.Item(6, 28).Value = 34
$ws.Cells.Item(6, 29).Value = 7.5
$ws.Cells.Item(6, 30).Value = 6.5
$ws.Cells.Item(6, 31).Value = 19
$ws.Cells.Item(6, 32).Value = 67
$ws.Cells.Item(6, 33).Value = 351
$ws.Cells.Item(6, 34).Value = 12
$ws.Cells.Item(6, 35).Value = 26
$ws.Cells.Item(6, 36).Value = 19
$ws.Cells.Item(6, 37).Value = 51
$ws.Cells.Item(6, 38).Value = 41
$ws.Cells.Item(6, 39).Value = 51
$ws.Cells.Item(6, 40).Value = 3.5
$ws.Cells.Item(6, 41).Value = 9.5
$ws.Cells.Item(6, 42).Value = 23
$ws.Cells.Item(6, 43).Value = 34
$ws.Cells.Item(6, 44).Value = 51
$ws.Cells.Item(6, 45).Value = 201
$ws.Cells.Item(6, 46).Value = 2.5
$ws.Cells.Item(6, 47).Value = 9.5
$ws.Cells.Item(6, 48).Value = 67
$ws.Cells.Item(6, 49).Value = 6.5
$ws.Cells.Item(6, 50).Value = 29
$ws.Cells.Item(6, 51).Value = 41
$ws.Cells.Item(6, 52).Value = 126
$ws.Cells.Item(6, 53).Value = 151
$ws.Cells.Item(6, 54).Value = 351
$ws.Cells.Item(6, 55).Value = 51
$ws.Cells.Item(6, 56).Value = 51

# --- Apply odds updates to row 7 (previously row 6, Carabobo vs Monagas) ---
$ws.Cells.Item(7, 7).Value = 1.75
$ws.Cells.Item(7, 8).Value = 3.3
$ws.Cells.Item(7, 9).Value = 4.6
$ws.Cells.Item(7, 14).Value = 6.6
$ws.Cells.Item(7, 15).Value = 1.35
$ws.Cells.Item(7, 22).Value = 1.72
$ws.Cells.Item(7, 23).Value = 6.1
$ws.Cells.Item(7, 24).Value = 7.7
$ws.Cells.Item(7, 26).Value = 14
$ws.Cells.Item(7, 28).Value = 30
$ws.Cells.Item(7, 29).Value = 8.25
$ws.Cells.Item(7, 30).Value = 6.5
$ws.Cells.Item(7, 34).Value = 10.75
$ws.Cells.Item(7, 35).Value = 25
$ws.Cells.Item(7, 37).Value = 90
$ws.Cells.Item(7, 38).Value = 55
$ws.Cells.Item(7, 44).Value = 60
$ws.Cells.Item(7, 47).Value = 7.4
$ws.Cells.Item(7, 49).Value = 6.3

